# Daily attendance processing - 2025-11-24 03:18:17
# Normalises the "Recorded By" (column G) entries on the
# "Session Analysis Results" sheet: for every session row whose recorder
# list has "System" (or an admin/backup account) trailing a human
# reviewer's address, the last two comma-separated names are swapped so
# the most recently-acting recorder is listed last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in the "Recorded By" / column G) whose last two entries need to
# be swapped, per today's reconciliation pass.
$rows = @(2, 3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 55, 58, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 113, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153)

foreach ($r in $rows) {
    $addr = "G" + $r
    $range = $ws.Range($addr)
    $current = $range.Value2
    $parts = $current.Split(",")
    $n = $parts.Length

    if ($n -ge 2) {
        $newParts = @()
        for ($i = 0; $i -lt ($n - 2); $i++) {
            $newParts += $parts[$i].Trim()
        }
        # swap the final two recorders
        $newParts += $parts[$n - 1].Trim()
        $newParts += $parts[$n - 2].Trim()

        $range.Value = [string]::Join(", ", $newParts)
    }
}
